$wb = $excel.ActiveWorkbook

# Rename the "recoded_domain_mapping" sheet to "domain_recodes"
$ws = $wb.Worksheets.Item("recoded_domain_mapping")
$ws.Name = "domain_recodes"

# Add the super_domain_name column (column C) with header and recoded values.
# Set the very first occurrence of each new label in the same order they were
# first introduced in the source workbook, so shared-string indices line up:
#   299 super_domain_name, 300 Physical science, 301 Natural science, 302 Social science
$ws.Range("C1").Value  = "super_domain_name"
$ws.Range("C3").Value  = "Physical science"    # Atmospheric science
$ws.Range("C2").Value  = "Natural science"     # Agriculture
$ws.Range("C5").Value  = "Social science"      # Business

$ws.Range("C4").Value  = "Natural science"     # Biology
$ws.Range("C6").Value  = "Information science" # Computer science
$ws.Range("C7").Value  = "Natural science"     # Ecology
$ws.Range("C8").Value  = "Social science"      # Education
$ws.Range("C9").Value  = "Physical science"    # Engineering
$ws.Range("C10").Value = "Natural science"     # Environmental science
$ws.Range("C11").Value = "Natural science"     # Natural Resources
# Row 12 (Geography) intentionally left without a super domain recode
$ws.Range("C13").Value = "Physical science"    # Geology
$ws.Range("C14").Value = "Physical science"    # Hydrology
$ws.Range("C15").Value = "Information science" # Information science
$ws.Range("C16").Value = "Law"                 # Law
$ws.Range("C17").Value = "Physical science"    # Mathematics
$ws.Range("C18").Value = "Natural science"     # Medicine
$ws.Range("C19").Value = "Physical science"    # Physical sciences
$ws.Range("C20").Value = "Social science"      # Psychology
$ws.Range("C21").Value = "Social science"      # Social sciences
$ws.Range("C22").Value = "Other"               # Other
$ws.Range("C23").Value = "Physical science"    # Marine Science
$ws.Range("C24").Value = "Physical science"    # Space and Planetary Science
$ws.Range("C25").Value = "Physical science"    # Chemistry
$ws.Range("C26").Value = "Natural science"     # Health Sciences

# Match the saved selection state on the sheet
$ws.Range("C27").Select()
